# Insert 4 new rows before row 10 (pushing existing rows 10:131 down to 14:135)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:13").Insert()

# Copy the constant "dimension" columns (A,B,C,E,F,G,H,I,J,K,R) from the row
# that now sits right below the newly inserted block (row 14), since these
# values are identical for every record in this sheet.
$constCols = @("A","B","C","E","F","G","H","I","J","K","R")
foreach ($col in $constCols) {
    $srcVal = $ws.Range($col + "14").Value()
    for ($r = 10; $r -le 13; $r++) {
        $ws.Range($col + $r).Value = $srcVal
    }
}

# New data block (week of 2021-12-17 / serial 44547)
$newRows = @{
    10 = @{ L = "Especial"; M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/caja 10 unidades"; S = 1750; T = 10 }
    11 = @{ L = "Primera";  M = 270; N = 17000; O = 18000; P = 17500; Q = "`$/caja 12 unidades"; S = 1458; T = 12 }
    12 = @{ L = "Segunda";  M = 270; N = 17000; O = 18000; P = 17500; Q = "`$/caja 14 unidades"; S = 1250; T = 14 }
    13 = @{ L = "Tercera";  M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/caja 16 unidades"; S = 1094; T = 16 }
}

for ($r = 10; $r -le 13; $r++) {
    $ws.Range("D" + $r).Value = 44547
    $row = $newRows[$r]
    $ws.Range("L" + $r).Value = $row.L
    $ws.Range("M" + $r).Value = $row.M
    $ws.Range("N" + $r).Value = $row.N
    $ws.Range("O" + $r).Value = $row.O
    $ws.Range("P" + $r).Value = $row.P
    $ws.Range("Q" + $r).Value = $row.Q
    $ws.Range("S" + $r).Value = $row.S
    $ws.Range("T" + $r).Value = $row.T
}
